$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.665.90'
$ws.Range('E2').Value = '  +7.00%  '
$ws.Range('D3').Value = '1.734.90'
$ws.Range('E3').Value = '  +3.50%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.87'
$ws.Range('E5').Value = '  +1.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9965'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3738'
$ws.Range('E7').Value = '  +2.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3394'
$ws.Range('E8').Value = '  +4.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.19'
$ws.Range('E9').Value = '  +2.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.185'
$ws.Range('E10').Value = '  +3.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07458'
$ws.Range('E11').Value = '  +5.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9956'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.403'
$ws.Range('E13').Value = '  +4.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.40'
$ws.Range('E14').Value = '  +3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.063'
$ws.Range('E15').Value = '  +6.39%  '
$ws.Range('D16').Value = '1.735.07'
$ws.Range('E16').Value = '  +3.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001075'
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06721'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.56'
$ws.Range('E19').Value = '  +4.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9974'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.65'
$ws.Range('E21').Value = '  +4.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.216'
$ws.Range('E22').Value = '  +4.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.77'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').Value = '26.757.80'
$ws.Range('E24').Value = '  +7.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.451'
$ws.Range('E25').Value = '  -0.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.471'
$ws.Range('E26').Value = '  +23.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.435'
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.66'
$ws.Range('E28').Value = '  +2.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.55'
$ws.Range('E29').Value = '  +3.83%  '
$ws.Range('D30').Value = '1.928.54'
$ws.Range('E30').Value = '  +3.58%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '132.24'
$ws.Range('E31').Value = '  +4.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.091'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.039'
$ws.Range('E33').Value = '  +4.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08599'
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.698'
$ws.Range('E35').Value = '  +2.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '12.86'
$ws.Range('E36').Value = '  +4.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.415'
$ws.Range('E37').Value = '  +4.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02352'
$ws.Range('E38').Value = '  +3.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2176'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06255'
$ws.Range('E40').Value = '  +3.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.442'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.224'
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6274'
$ws.Range('E43').Value = '  +5.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.27'
$ws.Range('E44').Value = '  +4.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9964'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.926'
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6080'
$ws.Range('E47').Value = '  +6.06%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '128.98'
$ws.Range('E48').Value = '  +2.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.070'
$ws.Range('E49').Value = '  +5.21%  '
$ws.Range('E50').Value = '  +2.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '77.64'
$ws.Range('E51').Value = '  +3.65%  '
